$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These price cells get new values that look like plain numbers. The
# source data stores them as text (so formats such as trailing zeros,
# e.g. "1.00" or "5.80", are preserved), so force a text format before
# writing them and restore the default cell style afterwards.
$textCells = @("D5", "D6", "D7", "D12", "D13", "D15", "D21", "D22", "D24", "D25", "D27", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D46", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.162.65"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "3.505.01"
$ws.Range("E3").Value = "  -4.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "581.25"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "174.03"
$ws.Range("E6").Value = "  -3.87%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "3.499.47"
$ws.Range("E8").Value = "  -4.85%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E11").Value = "  +6.15%  "
$ws.Range("D12").Value = "0.597"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "47.14"
$ws.Range("E13").Value = "  -5.93%  "
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "676.61"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "4.069.39"
$ws.Range("E16").Value = "  -5.11%  "
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "69.142.96"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("D19").Value = "3.510.67"
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "17.47"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").Value = "11.19"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("E23").Value = "  -4.42%  "
$ws.Range("D24").Value = "16.13"
$ws.Range("E24").Value = "  -9.62%  "
$ws.Range("D25").Value = "97.92"
$ws.Range("E25").Value = "  -5.94%  "
$ws.Range("E26").Value = "  -4.62%  "
$ws.Range("D27").Value = "5.83"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -6.77%  "
$ws.Range("D30").Value = "9.44"
$ws.Range("E30").Value = "  -7.36%  "
$ws.Range("D31").Value = "32.87"
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("E32").Value = "  -5.93%  "
$ws.Range("D33").Value = "3.20"
$ws.Range("E33").Value = "  -8.13%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "7.28"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "1.35"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").Value = "595.42"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("D37").Value = "3.60"
$ws.Range("E37").Value = "  -15.04%  "
$ws.Range("D38").Value = "10.90"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").Value = "57.32"
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -6.18%  "
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("D45").Value = "3.419.48"
$ws.Range("E45").Value = "  -10.15%  "
$ws.Range("D46").Value = "33.43"
$ws.Range("E46").Value = "  -6.21%  "
$ws.Range("E47").Value = "  -9.01%  "
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  -7.28%  "
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "5.80"
$ws.Range("E51").Value = "  +18.46%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
